# Weekly update: insert 3 new "Cultivar IV Región" rows for Chirimoya
# (market date 2021-10-25 / serial 44483, "Provincia del Elquí") at the
# top of the price block, pushing the existing history down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above row 20 (row 19 stays put, everything from the
# old row 20 onward shifts down to row 23 onward).
$ws.Rows("20:22").Insert()

# Common/static columns shared by every data row in this block.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107002
$categoria   = "Chirimoya"
$variedad    = "Cultivar IV Región"
$unidad      = "`$/kilo (en caja de 15 kilos)"
$origen      = "Provincia del Elquí"
$fecha       = 44483

$rows = @(
    @{ Row = 20; Calidad = "Especial"; Volumen = 240; PMin = 2200; PMax = 2300; PProm = 2250; PKg = 2250 },
    @{ Row = 21; Calidad = "Primera";  Volumen = 360; PMin = 1900; PMax = 2000; PProm = 1950; PKg = 1950 },
    @{ Row = 22; Calidad = "Segunda";  Volumen = 300; PMin = 1400; PMax = 1500; PProm = 1450; PKg = 1450 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = 1
}
